$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Range("C1").End(4).Row  ## xlDown = 4, walk from header to last filled cell in column C

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46060) {
        $cell.Value2 = 46061
    }
}
